# Add create stg_examples_images script and schema in spreadsheet
#
# This reproduces the same structural edit applied to Table10 / the
# stg_examples_images sheet as was previously done for stg_coins (Table7)
# and stg_examples (Table8): add "is_null" and "default" columns, and a
# calculated "sql_code" column. The original 5-column, 6-row table grows
# into an 8-column, 7-row table (header moved from row 1 to row 2; the
# original row-1 header cells are left behind as plain, non-table values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stg_examples_images")
$lo = $ws.ListObjects.Item("Table10")

# Grow the table down-and-right: new top-left is A2 (old header row becomes
# row 2), and it now spans 8 columns (A:H) and 6 data rows + header (7 total).
$lo.Resize($ws.Range("A2:H7"))

# --- Row 1: leftover header-shaped row (outside the table range) ---
$ws.Range("A1").Value = "table_name"
$ws.Range("B1").Value = "field_name"
$ws.Range("C1").Value = "data_type"
$ws.Range("D1").Value = "is_null"
$ws.Range("E1").Value = "default"
$ws.Range("F1").Value = "primary_key"
$ws.Range("G1").Value = "source"
$ws.Range("H1").Value = "sql_code"

# --- Row 2: the table's real header row ---
$ws.Range("A2").Value = "table_name"
$ws.Range("B2").Value = "field_name"
$ws.Range("C2").Value = "data_type"
$ws.Range("D2").Value = "is_null"
$ws.Range("E2").Value = "default"
$ws.Range("F2").Value = "primary_key"
$ws.Range("G2").Value = "source"
$ws.Range("H2").Value = "sql_code"

# --- Row 3: examples_images_id / SERIAL / NOT NULL / PK / Database generated
$ws.Range("A3").Value = "stg_examples_images"
$ws.Range("B3").Value = "examples_images_id"
$ws.Range("C3").Value = "SERIAL"
$ws.Range("D3").Value = "NOT NULL"
$ws.Range("F3").Value = $true
$ws.Range("G3").Value = "Database generated"

# --- Row 4: stg_examples_id / INTEGER / NOT NULL / not PK / FK source
$ws.Range("A4").Value = "stg_examples_images"
$ws.Range("B4").Value = "stg_examples_id"
$ws.Range("C4").Value = "INTEGER"
$ws.Range("D4").Value = "NOT NULL"
$ws.Range("F4").Value = $false
$ws.Range("G4").Value = "table: stg_examples, field: examples_id"

# --- Row 5: image_type / VARCHAR / NOT NULL / not PK / enum note
$ws.Range("A5").Value = "stg_examples_images"
$ws.Range("B5").Value = "image_type"
$ws.Range("C5").Value = "VARCHAR"
$ws.Range("D5").Value = "NOT NULL"
$ws.Range("F5").Value = $false
$ws.Range("G5").Value = "One of: obverse, reverse, both sides, unknown"

# --- Row 6: link / VARCHAR / NOT NULL / not PK / Scraping URI page
$ws.Range("A6").Value = "stg_examples_images"
$ws.Range("B6").Value = "link"
$ws.Range("C6").Value = "VARCHAR"
$ws.Range("D6").Value = "NOT NULL"
$ws.Range("F6").Value = $false
$ws.Range("G6").Value = "Scraping URI page"

# --- Row 7: ts / TIMESTAMP / NOT NULL / default CURRENT_TIMESTAMP / not PK
$ws.Range("A7").Value = "stg_examples_images"
$ws.Range("B7").Value = "ts"
$ws.Range("C7").Value = "TIMESTAMP"
$ws.Range("D7").Value = "NOT NULL"
$ws.Range("E7").Value = "CURRENT_TIMESTAMP"
$ws.Range("F7").Value = $false
$ws.Range("G7").Value = "Database generated"

# --- sql_code calculated column (same formula pattern as Table7 / Table8) ---
$ws.Range("H3").Formula = '=_xlfn.CONCAT(Table10[[#This Row],[field_name]], " ", Table10[[#This Row],[data_type]], " ", Table10[[#This Row],[is_null]], IF(LEN(Table10[[#This Row],[default]])=0,""," DEFAULT "&Table10[[#This Row],[default]]), ",")'
$ws.Range("H4").Formula = '=_xlfn.CONCAT(Table10[[#This Row],[field_name]], " ", Table10[[#This Row],[data_type]], " ", Table10[[#This Row],[is_null]], IF(LEN(Table10[[#This Row],[default]])=0,""," DEFAULT "&Table10[[#This Row],[default]]), ",")'
$ws.Range("H5").Formula = '=_xlfn.CONCAT(Table10[[#This Row],[field_name]], " ", Table10[[#This Row],[data_type]], " ", Table10[[#This Row],[is_null]], IF(LEN(Table10[[#This Row],[default]])=0,""," DEFAULT "&Table10[[#This Row],[default]]), ",")'
$ws.Range("H6").Formula = '=_xlfn.CONCAT(Table10[[#This Row],[field_name]], " ", Table10[[#This Row],[data_type]], " ", Table10[[#This Row],[is_null]], IF(LEN(Table10[[#This Row],[default]])=0,""," DEFAULT "&Table10[[#This Row],[default]]), ",")'
$ws.Range("H7").Formula = '=_xlfn.CONCAT(Table10[[#This Row],[field_name]], " ", Table10[[#This Row],[data_type]], " ", Table10[[#This Row],[is_null]], IF(LEN(Table10[[#This Row],[default]])=0,""," DEFAULT "&Table10[[#This Row],[default]]), ",")'

# Column widths: widen C:F, add the wide H (sql_code) column like on the
# sibling stg_coins / stg_examples sheets.
$ws.Columns.Item(8).ColumnWidth = 51.5

# --- View state: this sheet becomes the active tab, scrolled right so the
# new columns are visible, with the new formula cell selected.
$ws.Activate()
$ws.Range("F8").Select()
$excel.ActiveWindow.ScrollColumn = 3

# The sibling stg_examples sheet loses tabSelected and scrolls back to show
# its own sql_code column, matching the state captured in the diff.
$wsExamples = $wb.Worksheets.Item("stg_examples")
$wsExamples.Activate()
$wsExamples.Range("H2").Select()
$excel.ActiveWindow.ScrollColumn = 3

# Re-activate stg_examples_images last so it is the workbook's active sheet.
$ws.Activate()
